$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 36-38: Coin/Link values are reordered (shifted),
# plus Price/Volume values updated. Price values are numeric-looking
# text, so a leading apostrophe is used to force them to stay text
# (matches how Excel itself keeps typed numeric-looking text as text).
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'0.888"
$ws.Range("E36").Value = "  +7.83%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.14"
$ws.Range("E37").Value = "  +1.98%  "

$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.879"
$ws.Range("E38").Value = "  +0.40%  "

# Rows 47-48: Coin/Link values are swapped,
# plus Price/Volume values updated.
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'19.20"
$ws.Range("E47").Value = "  +2.55%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.63"
$ws.Range("E48").Value = "  +0.18%  "

# Remaining rows: update Price (D) and/or Volume (E) values only.
$ws.Range("D2").Value = "59.316.73"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "2.586.31"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'569.45"
$ws.Range("E5").Value = "  +2.89%  "

$ws.Range("D6").Value = "'143.63"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").Value = "2.596.26"
$ws.Range("E9").Value = "  -0.57%  "

$ws.Range("D10").Value = "'6.67"
$ws.Range("E10").Value = "  -1.71%  "

$ws.Range("E11").Value = "  +3.17%  "

$ws.Range("E12").Value = "  +9.10%  "

$ws.Range("D13").Value = "'0.344"
$ws.Range("E13").Value = "  +2.23%  "

$ws.Range("D14").Value = "3.041.81"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").Value = "59.343.49"
$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("D16").Value = "'22.56"
$ws.Range("E16").Value = "  +7.95%  "

$ws.Range("E17").Value = "  +4.16%  "

$ws.Range("D18").Value = "2.590.08"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").Value = "'4.52"
$ws.Range("E19").Value = "  +1.57%  "

$ws.Range("D20").Value = "'336.54"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").Value = "'10.22"
$ws.Range("E21").Value = "  +1.48%  "

$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'64.19"
$ws.Range("E24").Value = "  -3.42%  "

$ws.Range("D25").Value = "'0.453"
$ws.Range("E25").Value = "  +5.92%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("D28").Value = "'7.27"
$ws.Range("E28").Value = "  +1.90%  "

$ws.Range("E29").Value = "  +3.71%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").Value = "'6.08"
$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("D33").Value = "'156.99"
$ws.Range("E33").Value = "  +2.62%  "

$ws.Range("D34").Value = "'19.06"
$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("D35").Value = "'4.04"
$ws.Range("E35").Value = "  +3.04%  "

$ws.Range("E39").Value = "  +2.83%  "

$ws.Range("D40").Value = "'36.82"
$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("D41").Value = "'295.19"
$ws.Range("E41").Value = "  +4.11%  "

$ws.Range("D42").Value = "'3.67"
$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "'0.0975"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").Value = "'0.598"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("D46").Value = "'0.0537"
$ws.Range("E46").Value = "  +0.87%  "

$ws.Range("D49").Value = "'124.66"
$ws.Range("E49").Value = "  +4.95%  "

$ws.Range("E50").Value = "  +2.47%  "

$ws.Range("D51").Value = "'18.56"
$ws.Range("E51").Value = "  +4.09%  "
